# This workbook tracks daily "Bruselas (repollito)" price records.
# A new weekly record needs to be inserted as row 41 (pushing the
# existing rows 41-56 down to 42-57), matching the new sample taken on
# 2022-08-11 (serial date 44784).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 41; this shifts existing rows 41-56 down
# to 42-57 and expands the sheet dimension automatically.
$ws.Rows.Item(41).Insert()

# Columns A, B, C, E, F, G, H, I, N, O, Q, R are constant across every
# record in this table, so copy them from the row directly below
# (which now holds what used to be row 41).
$ws.Cells.Item(41, 1).Value2  = $ws.Cells.Item(42, 1).Value2   # A Mercado ID
$ws.Cells.Item(41, 2).Value2  = $ws.Cells.Item(42, 2).Value2   # B Mercado
$ws.Cells.Item(41, 3).Value2  = $ws.Cells.Item(42, 3).Value2   # C Region
$ws.Cells.Item(41, 4).Value2  = 44784                          # D Fecha
$ws.Cells.Item(41, 5).Value2  = $ws.Cells.Item(42, 5).Value2   # E Codreg
$ws.Cells.Item(41, 6).Value2  = $ws.Cells.Item(42, 6).Value2   # F Categoria ID
$ws.Cells.Item(41, 7).Value2  = $ws.Cells.Item(42, 7).Value2   # G Categoria
$ws.Cells.Item(41, 8).Value2  = $ws.Cells.Item(42, 8).Value2   # H Variedad
$ws.Cells.Item(41, 9).Value2  = $ws.Cells.Item(42, 9).Value2   # I Calidad
$ws.Cells.Item(41, 10).Value2 = 28                              # J Volumen
$ws.Cells.Item(41, 11).Value2 = 20000                           # K Precio minimo
$ws.Cells.Item(41, 12).Value2 = 21000                           # L Precio maximo
$ws.Cells.Item(41, 13).Value2 = 20357                           # M Precio promedio ponderado
$ws.Cells.Item(41, 14).Value2 = $ws.Cells.Item(42, 14).Value2  # N Unidad de comercializacion
$ws.Cells.Item(41, 15).Value2 = $ws.Cells.Item(42, 15).Value2  # O Origen
$ws.Cells.Item(41, 16).Value2 = 1357                            # P Precio $/Kg
$ws.Cells.Item(41, 17).Value2 = $ws.Cells.Item(42, 17).Value2  # Q Kg o Unidades
$ws.Cells.Item(41, 18).Value2 = $ws.Cells.Item(42, 18).Value2  # R Clasificacion
